$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Tabela1")

# New daily COVID-19 data rows appended to the bottom of the table.
$newData = @(
    @(43976, 75770, 754, 1469, 0, 9, 2, 6, 108, 1),
    @(43977, 76579, 809, 1471, 2, 8, 2, 2, 108, 0),
    @(43978, 77210, 631, 1473, 2, 7, 2, 1, 108, 0)
)

foreach ($rowVals in $newData) {
    # Format the new row the same way as the table's current last row, then
    # let the table grow by one row (keeps the autofilter / dimension / ref
    # in sync automatically).
    $lastRowRange = $tbl.ListRows.Item($tbl.ListRows.Count).Range
    $lastRowRange.Copy()

    $newRow = $tbl.ListRows.Add()
    $destRow = $newRow.Range
    $destRow.PasteSpecial(-4122)

    for ($col = 1; $col -le 10; $col++) {
        $destRow.Cells.Item(1, $col).Value2 = $rowVals[$col - 1]
    }
}

$excel.CutCopyMode = 0

$lastRow = $tbl.Range.Row + $tbl.Range.Rows.Count - 1
$ws.Range("A" + $lastRow + ":J" + $lastRow).Select()
